$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Table 1: Market Risk Premium sensitivity (K3:L8) ---
$ws.Range("K3").Value = "Market Risk Premium"
$ws.Range("L3").Value = "WACC"

$mrpK = @(0.03, 0.04, 0.05, 0.06, 0.07)
$mrpL = @(0.069708, 0.081444, 0.09318, 0.104916, 0.116652)

for ($i = 0; $i -lt $mrpK.Length; $i++) {
    $row = 4 + $i
    $kCell = $ws.Cells.Item($row, 11)
    $lCell = $ws.Cells.Item($row, 12)
    $kCell.Value = $mrpK[$i]
    $lCell.Value = $mrpL[$i]
    $kCell.NumberFormat = "0%"
    $lCell.NumberFormat = "0.00%"
}

# --- Table 2: Tax Rate sensitivity (K10:L21) ---
$ws.Range("K10").Value = "Tax Rate"

$taxK = @(0.35, 0.36, 0.37, 0.38, 0.39, 0.4, 0.41, 0.42, 0.43, 0.44, 0.45)
$taxL = @(0.093595, 0.093512, 0.093429, 0.093346, 0.093263, 0.09318, 0.093097, 0.093014, 0.092931, 0.092848, 0.092765)

for ($i = 0; $i -lt $taxK.Length; $i++) {
    $row = 11 + $i
    $kCell = $ws.Cells.Item($row, 11)
    $lCell = $ws.Cells.Item($row, 12)
    $kCell.Value = $taxK[$i]
    $lCell.Value = $taxL[$i]
    $kCell.NumberFormat = "0%"
    $lCell.NumberFormat = "0.00%"
}

# Match the author's final selection state
[void]$ws.Range("M10").Select()
